$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 5 (old row 4 will be duplicated there); ---
# --- old rows 5,6,7 shift down to 6,7,8 automatically.                     ---
$ws.Rows(5).Insert()

# --- 2. Duplicate the original row 4 ("...in the last week" question) into the ---
# ---    newly inserted row 5, preserving its values/number formats/styles.    ---
$ws.Range("A4:E4").Copy($ws.Range("A5:E5"))
$ws.Rows(5).RowHeight = 409.6

# --- 3. Hyperlinks: the RAND URL hyperlink that lived on D4 now belongs to the ---
# ---    duplicated row (D5); D4 will get its own, newly created hyperlink.    ---
$url = "https://www.rand.org/health-care/surveys_tools/mos/social-support.html"
$ws.Range("D4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D5"), $url) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), $url) | Out-Null

# --- 4. Replace A4 with the new question text that asks about the period      ---
# ---    "prior to the COVID-19 pandemic" instead of "in the last week".       ---
$ws.Range("A4").Value = 'People sometimes look to others for companionship, assistance, or other types of support. How often is each of the following kinds of support available to you if you need it? Please choose one number that describes how often each kind of support was available to in a typical week prior to the COVID-19 pandemic. 
•	Someone you can count on to listen to you when you need to talk
•	Someone to give you information to help you understand a situation
•	Someone to give you good advice about a crisis
•	Someone to confide in or talk to about yourself or your problems
•	Someone whose advice you really want
•	Someone to share your most private worries and fears with
•	Someone to turn to for suggestions about how to deal with a personal problem
•	Someone who understands your problems
•	Someone to help you if you were confined to bed
•	Someone to take you to the doctor if you needed it
•	Someone to prepare your meals if you were unable to do it yourself
•	Someone to help with daily chores if you were sick
•	Someone who shows you love and affection
•	Someone to love and make you feel wanted
•	Someone who hugs you
•	Someone to have a good time with
•	Someone to get together with for relaxation
•	Someone to do something enjoyable with
•	Someone to do things with to help you get your mind off things'
$a4 = $ws.Range("A4")
$a4.Characters(265, 48).Font.Bold = $true
$a4.Characters(265, 48).Font.Size = 12
$a4.Characters(313, 991).Font.Size = 12

# --- 5. Update the sheet view: scroll down a bit and select B4 (matches the   ---
# ---    author's saved cursor position after editing this area).             ---
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("B4").Select()

Write-Host "Edit complete"
